$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.384.23'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.13%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.879.51'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.03%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7168'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.73%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '243.73'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.63%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("E8").Value = '  -1.43%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3145'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.62%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.92'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.28%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08084'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -3.66%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.878.76'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.25%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '94.74'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +3.64%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.220'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.61%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.7083'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.53%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.382'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.41%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008422'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.44%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.379.63'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.16%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '252.45'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +4.71%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.34'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.67%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.132.02'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.02%  '

$ws.Range("E22").Value = '  +0.13%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.682'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.50%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.00%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1580'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.94%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.063'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.01%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.65'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.78%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.98'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +2.19%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.511'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.31%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.419'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.18%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.318'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.76%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.230'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.93%  '

$ws.Range("E33").Value = '  -1.22%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.941'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.44%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7583'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.84%  '

$ws.Range("E36").Value = '  -0.27%  '

$ws.Range("E37").Value = '  +0.27%  '

$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.290.84'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.16%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01883'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.33%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.766'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.07%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.406'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.69%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9063'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +1.59%  '

$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '74.08'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.03%  '

$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '111.50'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.46%  '

$ws.Range("E45").Value = '  +0.03%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000129'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.55%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.027.19'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.25%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.808'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.31%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.5206'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.02%  '

$ws.Range("E50").Value = '  +0.58%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4348'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.42%  '
